# The sheet "Đơn sale chính" (first/active sheet) gets a new, blank second
# row appended beneath its header row (dimension grows from A1:T1 to
# A1:T2). Several of the numeric columns in that new row carry an explicit
# 0 value while the rest of the row stays empty.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Touch every cell A2:T2 (with a no-op "no border" formatting pass) so the
# whole row materializes in the sheet, mirroring the header row's span,
# before filling in the numeric columns that actually carry a value.
$ws.Range("A2:T2").Borders.LineStyle = -4142

$ws.Range("B2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
